$wb = $excel.ActiveWorkbook

# Sheets that contain the duplicated event listing data: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 5107
    $ws.Range("F4").Value = 891
}
